{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1. Title (Heading1 paragraph)\nparagraphs.items[0].insertText(\n  \"Review 197: [Short] LLM4Decompile: Decompiling Binary Code with Large Language Models\",\n  \"Replace\"\n);\n\n// 2. Bold \"Paper: ...\" line\nparagraphs.items[1].insertText(\n  \"Paper: https://arxiv.org/abs/2403.05286v3\",\n  \"Replace\"\n);\n\n// 3. Plain arxiv abstract link paragraph\nparagraphs.items[3].insertText(\n  \"https://arxiv.org/abs/2403.05286\",\n  \"Replace\"\n);\n\n// 4. Opening Hebrew summary paragraph -> replaced with new LLM4Decompile summary\nparagraphs.items[5].insertText(\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 LLM4Decompile, \u05de\u05e9\u05e4\u05d7\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 LLM \u05dc\u05d3\u05d9\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d1\u05d2\u05d9\u05e9\u05d4 \u05e4\u05ea\u05d5\u05d7\u05d4 \u05e9\u05e0\u05e2\u05d9\u05dd \u05de-1B \u05e2\u05d3 33B \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd. \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05e2\u05dc 4 \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05d5\u05d3 \u05de\u05e7\u05d5\u05e8 \u05d1\u05e9\u05e4\u05ea C \u05d5\u05e7\u05d5\u05d3 \u05d0\u05e1\u05de\u05d1\u05dc\u05d9 \u05de\u05ea\u05d0\u05d9\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d2\u05dd \u05de\u05e6\u05d9\u05d2\u05d9\u05dd \u05d0\u05ea Decompile-Eval, \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05d4\u05e2\u05e8\u05db\u05ea \u05d3\u05d9\u05d5\u05e7 \u05d4\u05d3\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d4\u05de\u05d1\u05d5\u05e6\u05e2\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d5\u05d3\u05dc (\u05de\u05e7\u05de\u05e4\u05dc\u05d9\u05dd \u05de\u05d7\u05d3\u05e9 \u05d5\u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e7\u05d5\u05d3). \",\n  \"Replace\"\n);\n\n// 5. Remove the large block of old review paragraphs (indices 7 through 24 inclusive),\n//    going in reverse order so earlier indices stay valid as we delete.\nfor (let i = 24; i >= 7; i--) {\n  paragraphs.items[i].delete();\n}\n\n// 6. The final paragraph (formerly \"\u05d5\u05d6\u05d4 \u05d5\u05d6\u05d4\u2026\") now holds the new closing summary text\nparagraphs.items[25].insertText(\n  \"LLM4Decompile \u05de\u05e6\u05dc\u05d9\u05d7 \u05dc\u05d1\u05e6\u05e2 \u05d3\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d3\u05d5\u05d9\u05e7\u05ea 21% \u05de\u05e7\u05d5\u05d3 \u05d4\u05d0\u05e1\u05de\u05d1\u05dc\u05d9, \u05e2\u05dd \u05e9\u05d9\u05e4\u05d5\u05e8 \u05e9\u05dc 50% \u05d1\u05d9\u05d7\u05e1 \u05dc-GPT-4. \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05dc\u05e7\u05d9\u05de\u05e4\u05d5\u05dc \u05d5\u05dc\u05d3\u05e7\u05de\u05e4\u05d5\u05dc \u05e9\u05dc \u05e7\u05d5\u05d3 \u05e0\u05e8\u05d0\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05d7\u05d5\u05dd \u05de\u05d7\u05e7\u05e8 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05e1\u05d5\u05db\u05e0\u05d9 AI \u05d7\u05e1\u05d9\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e0\u05d2\u05d3 \u05d4\u05ea\u05e7\u05e4\u05d5\u05ea \u05d0\u05d3\u05d5\u05d5\u05e8\u05e1\u05e8\u05d9\u05d5\u05ea \u05e9\u05d9\u05e4\u05e2\u05dc\u05d5 \u05d1\u05e9\u05db\u05d1\u05d5\u05ea \u05d4\u05e2\u05de\u05d5\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc Software Stack.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Title (Heading1 paragraph)\n$d.Paragraphs.Item(1).Range.Text = \"Review 197: [Short] LLM4Decompile: Decompiling Binary Code with Large Language Models\"\n\n# 2. Bold \"Paper: ...\" line\n$d.Paragraphs.Item(2).Range.Text = \"Paper: https://arxiv.org/abs/2403.05286v3\"\n\n# 3. Plain arxiv abstract link paragraph\n$d.Paragraphs.Item(4).Range.Text = \"https://arxiv.org/abs/2403.05286\"\n\n# 4. Opening Hebrew summary paragraph -> replaced with new LLM4Decompile summary\n$d.Paragraphs.Item(6).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 LLM4Decompile, \u05de\u05e9\u05e4\u05d7\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 LLM \u05dc\u05d3\u05d9\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d1\u05d2\u05d9\u05e9\u05d4 \u05e4\u05ea\u05d5\u05d7\u05d4 \u05e9\u05e0\u05e2\u05d9\u05dd \u05de-1B \u05e2\u05d3 33B \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd. \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d0\u05dc\u05d5 \u05de\u05d0\u05d5\u05de\u05e0\u05d9\u05dd \u05e2\u05dc 4 \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc \u05e7\u05d5\u05d3 \u05de\u05e7\u05d5\u05e8 \u05d1\u05e9\u05e4\u05ea C \u05d5\u05e7\u05d5\u05d3 \u05d0\u05e1\u05de\u05d1\u05dc\u05d9 \u05de\u05ea\u05d0\u05d9\u05dd. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d2\u05dd \u05de\u05e6\u05d9\u05d2\u05d9\u05dd \u05d0\u05ea Decompile-Eval, \u05d4\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05dc\u05d4\u05e2\u05e8\u05db\u05ea \u05d3\u05d9\u05d5\u05e7 \u05d4\u05d3\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d4\u05de\u05d1\u05d5\u05e6\u05e2\u05ea \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d5\u05d3\u05dc (\u05de\u05e7\u05de\u05e4\u05dc\u05d9\u05dd \u05de\u05d7\u05d3\u05e9 \u05d5\u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05e0\u05dc\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4\u05e7\u05d5\u05d3). \"\n\n# 5. Remove the large block of old review paragraphs (paragraphs 8 through 25, 1-based),\n#    which sits between the empty paragraph after the summary and the final paragraph.\n$pStart = $d.Paragraphs.Item(8)\n$pEnd = $d.Paragraphs.Item(25)\n$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)\n$rng.Delete()\n\n# 6. The final paragraph (formerly \"\u05d5\u05d6\u05d4 \u05d5\u05d6\u05d4\u2026\") now holds the new closing summary text\n$d.Paragraphs.Item(8).Range.Text = \"LLM4Decompile \u05de\u05e6\u05dc\u05d9\u05d7 \u05dc\u05d1\u05e6\u05e2 \u05d3\u05e7\u05d5\u05de\u05e4\u05d9\u05dc\u05e6\u05d9\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d3\u05d5\u05d9\u05e7\u05ea 21% \u05de\u05e7\u05d5\u05d3 \u05d4\u05d0\u05e1\u05de\u05d1\u05dc\u05d9, \u05e2\u05dd \u05e9\u05d9\u05e4\u05d5\u05e8 \u05e9\u05dc 50% \u05d1\u05d9\u05d7\u05e1 \u05dc-GPT-4. \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05dc\u05e7\u05d9\u05de\u05e4\u05d5\u05dc \u05d5\u05dc\u05d3\u05e7\u05de\u05e4\u05d5\u05dc \u05e9\u05dc \u05e7\u05d5\u05d3 \u05e0\u05e8\u05d0\u05d9\u05dd \u05db\u05de\u05d5 \u05ea\u05d7\u05d5\u05dd \u05de\u05d7\u05e7\u05e8 \u05d7\u05e9\u05d5\u05d1 \u05d1\u05de\u05d9\u05d5\u05d7\u05d3 \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05e1\u05d5\u05db\u05e0\u05d9 AI \u05d7\u05e1\u05d9\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e0\u05d2\u05d3 \u05d4\u05ea\u05e7\u05e4\u05d5\u05ea \u05d0\u05d3\u05d5\u05d5\u05e8\u05e1\u05e8\u05d9\u05d5\u05ea \u05e9\u05d9\u05e4\u05e2\u05dc\u05d5 \u05d1\u05e9\u05db\u05d1\u05d5\u05ea \u05d4\u05e2\u05de\u05d5\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc Software Stack.\"\n"}
